# Small Tweaks / View doesn`t work yet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product id numbers in column F (rows 29-40) from 101..112 to 1..12
for ($i = 0; $i -lt 12; $i++) {
    $row = 29 + $i
    $ws.Range("F$row").Value = $i + 1
}

# Update the view: scroll/top-left cell and selected cell, and zoom to 100%
$ws.Activate()
$ws.Range("G25").Select()
$excel.ActiveWindow.ScrollRow = 26
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.Zoom = 100
